# Edit script: reproduces the commit
#   - table on slide 6 switches from the custom "Table_0" style to the
#     built-in table style {3568D5AE-2369-4A6B-817C-5A4DF298FEF8}
#   - the presentation's applied theme ("Integral") is swapped out for the
#     default "Office Theme" colour scheme

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style on slide 6 (the table graphic frame is the 2nd shape)
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{3568D5AE-2369-4A6B-817C-5A4DF298FEF8}")
    }
}

# ---------------------------------------------------------------------
# 2) Swap the design's colour scheme from "Integral" to "Office Theme"
# ---------------------------------------------------------------------
$cs = $p.SlideMaster.ColorScheme
$cs.Colors(1).RGB  = 0          # dk1       000000
$cs.Colors(2).RGB  = 16777215   # lt1       FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2       44546A
$cs.Colors(4).RGB  = 15132391   # lt2       E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1   5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2   ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3   A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4   FFC000
$cs.Colors(9).RGB  = 12874308   # accent5   4472C4
$cs.Colors(10).RGB = 4697456    # accent6   70AD47
$cs.Colors(11).RGB = 12673797   # hlink     0563C1
$cs.Colors(12).RGB = 7491477    # folHlink  954F72
